# Add a new worksheet "ODI Bowling Extra" (mirrors the existing
# "ODI Batting Extra" sheet) as the LAST sheet in the workbook, and
# populate it with the MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL
# data scraped for the extra bowling attributes.

$wb = $excel.ActiveWorkbook

# --- tidy up "ODI Batting Extra": rows whose batting stats were never
#     scraped (B:E only ever held placeholder blanks) shouldn't carry
#     empty cells around now that we're adding real scraped data elsewhere ---
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")
$battingExtra.Range("B4:E5").ClearContents()
$battingExtra.Range("B12:E13").ClearContents()
$battingExtra.Range("B17:E18").ClearContents()

# --- create the sheet after the current last tab so it lands at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Bowling Extra"

# --- header row (bold, centered, top-aligned, boxed -- same look as the
#     header rows on the other scraped sheets) ---
$headers = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
for ($c = 1; $c -le $headers.Length; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = $headers[$c - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# --- data rows ---
# Each tuple is (MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL); blank
# strings mean the scraper found no value for that match/column.
$data = @(
    @("3973", "", ""),
    @("4108", "", ""),
    @("4115", "0", ""),
    @("4117", "0", "10.00%"),
    @("4123", "1", ""),
    @("4125", "0", "20.00%"),
    @("4415", "0", ""),
    @("4419", "", ""),
    @("4421", "", ""),
    @("4423", "0", "30.00%"),
    @("4429", "1", "10.00%"),
    @("4430", "1", "10.00%"),
    @("4431", "0", ""),
    @("4483", "0", "10.00%"),
    @("4484", "0", ""),
    @("4601", "1", "20.00%"),
    @("4603", "0", ""),
    @("4644", "0", "10.00%"),
    @("4663", "", ""),
    @("4666", "", "")
)

$r = 2
foreach ($row in $data) {
    for ($c = 1; $c -le 3; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        # Force text formatting first so numeric-looking values ("0", "1")
        # and percentages ("10.00%") are stored as literal text, matching
        # the rest of the scraped sheets instead of being auto-coerced to
        # numbers/percentages by Excel.
        $cell.NumberFormat = "@"
        $cell.Value = $row[$c - 1]
    }
    $r++
}

$ws.Range("A1").Select() | Out-Null
